$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header labels (row 1) ---
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

# --- Row 2: area + totals ---
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("H2").Formula = "=SUM(G2:G11)"
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# --- Re-fill D3:D9 as one shared-formula range (keeps same formula text) ---
$ws.Range("D3:D9").Formula = "=(A3/100+(A4/100-A3/100)/2)"

# --- G3 standalone formula ---
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"

# --- G4:G15 as one shared-formula range ---
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

# --- Update selection to match the new working cells ---
$ws.Range("J2:K2").Select()
